$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 2: "Lancez 3 builds" -> wrap "builds" in spellStart/spellEnd proofErr ---
$p = $d.Paragraphs(2)
$xml = '<w:p ' + $wns + ' w:rsidR="00037294" w:rsidRDefault="003A367C" w:rsidP="00037294">' + `
  '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Lancez 3 </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r w:rsidR="00037294"><w:t>builds</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 3: "Lancez la carte de droite" -> "Lancez la carte de gauche" (split run) ---
$p = $d.Paragraphs(3)
$xml = '<w:p ' + $wns + ' w:rsidR="003A367C" w:rsidRDefault="003A367C" w:rsidP="003A367C">' + `
  '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Lancez la carte de </w:t></w:r>' + `
  '<w:r><w:t>gauche</w:t></w:r>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 4: "Sur une de fenêtre, appuyez sur espace lors de l'interface de connextion"
#                   -> "Sur une de fenêtre, lancez un serveur" (split run) ---
$p = $d.Paragraphs(4)
$xml = '<w:p ' + $wns + ' w:rsidR="00037294" w:rsidRDefault="00037294" w:rsidP="003A367C">' + `
  '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Sur une de fenêtre, </w:t></w:r>' + `
  '<w:r><w:t>lancez un serveur</w:t></w:r>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 6: "Appuyez sur play" -> split + proofErr around "play" ---
$p = $d.Paragraphs(6)
$xml = '<w:p ' + $wns + ' w:rsidR="003A367C" w:rsidRDefault="00037294" w:rsidP="003A367C">' + `
  '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Appuyez sur </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>play</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 7: "Pendant la phase de reflexion, ..." -> split + proofErr around "reflexion" ---
$p = $d.Paragraphs(7)
$xml = '<w:p ' + $wns + ' w:rsidR="00037294" w:rsidRDefault="00037294" w:rsidP="003A367C">' + `
  '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Pendant la phase de </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>reflexion</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>, construisez des tourelles (carrés blanc) et demandez à récolter des ressources (bouton mission)</w:t></w:r>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 11: "Vous pouvez mettre pause avec echap" -> split + proofErr around "echap" ---
$p = $d.Paragraphs(11)
$xml = '<w:p ' + $wns + ' w:rsidR="006C4DF5" w:rsidRDefault="006C4DF5" w:rsidP="006C4DF5">' + `
  '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Vous pouvez mettre pause avec </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>echap</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 12: "... paramètres dans le menu (echap)" -> split last run + proofErr around "echap" ---
$p = $d.Paragraphs(12)
$xml = '<w:p ' + $wns + ' w:rsidR="006C4DF5" w:rsidRDefault="00037294" w:rsidP="006C4DF5">' + `
  '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Toutes les options sont débloqués </w:t></w:r>' + `
  '<w:r w:rsidR="000D0285"><w:t>dans cette version</w:t></w:r>' + `
  '<w:r><w:t>' + [char]0x00A0 + ': vous pouvez customisez des</w:t></w:r>' + `
  '<w:r w:rsidR="006C4DF5"><w:t xml:space="preserve"> paramètres dans le menu (</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>echap</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>)</w:t></w:r>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 13: "Vous pouvez spécialiser une tourelle à partir du niveau 4" ---
#     Drop the bookmark here (it moves to the new final paragraph below).
$p = $d.Paragraphs(13)
$xml = '<w:p ' + $wns + ' w:rsidR="00EF5EF8" w:rsidRDefault="00EF5EF8" w:rsidP="006C4DF5">' + `
  '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Vous pouvez spécialiser une tourelle à partir du niveau 4</w:t></w:r>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Append a blank paragraph, then the new closing paragraph carrying the _GoBack bookmark ---
$p = $d.Paragraphs(13)
$p.Range.InsertParagraphAfter()
$blank = $d.Paragraphs(14)
$blank.Range.InsertXML('<w:p ' + $wns + '/>')

$blank = $d.Paragraphs(14)
$blank.Range.InsertParagraphAfter()
$final = $d.Paragraphs(15)
$xml = '<w:p ' + $wns + '>' + `
  '<w:r><w:t>En cas d’erreur de connexion au serveur, vous pouvez vous connecter avec le client et avec le server en connexion locale (cocher la coche connexion local et pour les clients, entrez l’adresse 127.0.0.1).</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
$final.Range.InsertXML($xml)
